# Generate Report for Handback
#
# Row 7 of both the "zh-cn" and "de-de" sheets describes the handback for
# ba987542-8998-49c3-842f-aca4e7047f50 (status "Ready for handoff"). This
# edit records the outcome of validating that handback: the handed-back
# file's commit is not the latest one available, so we populate the
# "Latest Target File" (I), "Latest Handback File" (J), "Latest Handback
# DateTime" (K) and "Error Detail" (P) columns for that row on both
# language sheets.

$wb = $excel.ActiveWorkbook
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsDeDe = $wb.Worksheets.Item("de-de")

$currentUrl = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/17a067cdf135f23a4872e70c58852f20eb56c067/e2e/ba987542-8998-49c3-842f-aca4e7047f50.md"
$latestUrl = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/b1a973334e91269aa9f65fc835d9499de29403b5/e2e/ba987542-8998-49c3-842f-aca4e7047f50.md"
$handbackDisplay = "ba987542-8998-49c3-842f-aca4e7047f50.md"
$errorDetail = "The version of handback file is not the latest, current: $currentUrl, latest: $latestUrl."

# --- zh-cn sheet, row 7 ---
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("I7"), $currentUrl, "", "", $handbackDisplay) | Out-Null
$wsZhCn.Range("J7").Value = "ba987542-8998-49c3-842f-aca4e7047f50.12f8173d68b52a6d2d4251a84ba4982118f5171d.zh-cn.xlf"
$wsZhCn.Range("K7").Value = "2016-08-20 04:57:55"
$wsZhCn.Range("P7").Value = $errorDetail

# --- de-de sheet, row 7 ---
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("I7"), $currentUrl, "", "", $handbackDisplay) | Out-Null
$wsDeDe.Range("J7").Value = "ba987542-8998-49c3-842f-aca4e7047f50.12f8173d68b52a6d2d4251a84ba4982118f5171d.de-de.xlf"
$wsDeDe.Range("K7").Value = "2016-08-20 04:58:07"
$wsDeDe.Range("P7").Value = $errorDetail
